# "Generate Report for Handoff"
#
# Regenerates the localization-status report: the zh-cn / de-de handoff
# status moves from "Handed back: in sync with en-US" to "Ready for
# handoff", and the associated report/handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = latest HO xliff generate date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 10:59:43"

# --- zh-cn sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 10:59:37"

# --- de-de sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 10:59:43"

# --- Column width refresh ---------------------------------------------
# With the shorter "Ready for handoff" text replacing the much longer
# "Handed back: in sync with en-US", the report generator re-sized the
# Status columns narrower on all three sheets.
$wsOverview.Range("E1:F1").ColumnWidth = 16.3
$wsZhCn.Range("C1").ColumnWidth = 16.3
$wsDeDe.Range("C1").ColumnWidth = 16.3
